$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) "json" list (column M): insert new entry "storeKeys(json,jsonpath,var)"
#    in its alphabetically-correct spot, just before "storeValue(...)" and
#    "storeValues(...)" which both shift down by one row (M16->M17->M18).
# ---------------------------------------------------------------------------
$ws.Range("M18").Value2 = $ws.Range("M17").Value()
$ws.Range("M17").Value2 = $ws.Range("M16").Value()
$ws.Range("M16").Value2 = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------------
# 2) "target" list (column A): remove the "text" category entry (was A25),
#    shifting web / webalert / webcookie / ws / ws.async / xml up by one row.
# ---------------------------------------------------------------------------
for ($r = 25; $r -le 30; $r++) {
    $ws.Range("A$r").Value2 = $ws.Range("A$($r + 1)").Value()
}
$ws.Range("A31").ClearContents()

# ---------------------------------------------------------------------------
# 3) Drop the "text" category's function column entirely (old column Y held
#    only spellCheck(var,profile,text)); everything from Z..AE shifts left
#    into Y..AD.
# ---------------------------------------------------------------------------
$ws.Columns.Item(25).Delete()

# ---------------------------------------------------------------------------
# 4) Keep the named ranges in sync with the new layout.
# ---------------------------------------------------------------------------
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
